# Update the two-digit multiplication equations in the document.
# Each old equation string is unique within the document, so a simple
# Find/Replace (wildcards off, match case on, whole-document range) is
# sufficient and unambiguous for every pair.

$d = $word.ActiveDocument

$pairs = @(
    @("69×78=5382", "91×81=7371"),
    @("65×54=3510", "66×79=5214"),
    @("57×58=3306", "67×56=3752"),
    @("66×77=5082", "57×68=3876"),
    @("50×72=3600", "59×49=2891"),
    @("11×83=913",  "30×34=1020"),
    @("64×79=5056", "54×28=1512"),
    @("53×14=742",  "50×71=3550"),
    @("42×65=2730", "14×59=826"),
    @("69×48=3312", "86×77=6622"),
    @("57×67=3819", "59×25=1475"),
    @("98×59=5782", "45×24=1080"),
    @("43×43=1849", "59×45=2655"),
    @("59×61=3599", "21×57=1197"),
    @("99×67=6633", "65×88=5720"),
    @("54×66=3564", "57×62=3534"),
    @("39×96=3744", "60×50=3000"),
    @("41×63=2583", "60×52=3120"),
    @("52×62=3224", "13×69=897"),
    @("68×65=4420", "31×49=1519"),
    @("68×43=2924", "63×68=4284"),
    @("13×91=1183", "59×29=1711"),
    @("82×82=6724", "19×50=950"),
    @("33×43=1419", "61×45=2745"),
    @("62×93=5766", "82×86=7052")
)

foreach ($pair in $pairs) {
    $oldText = $pair[0]
    $newText = $pair[1]

    $range = $d.Content
    $range.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)
}

$d.Save()
